# Update "gh-pages" output - refresh scraped counts and insert a newly
# discovered event ("广州·第六届淋唔到动漫嘉年华", 2024-06-22) into both the
# "展览" (Exhibition) sheet and the "全部类型" (All types) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (index 1) - bump a batch of "want to go" counters (col F)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$updates1 = @{
    3  = 148
    4  = 1779
    6  = 1061
    7  = 2208
    8  = 2125
    9  = 1111
    10 = 607
    15 = 44
    16 = 97
    18 = 1588
    20 = 727
    22 = 12258
    23 = 12310
    24 = 910
    27 = 38
    29 = 371
    30 = 1923
}

foreach ($row in $updates1.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates1[$row]
}

# Insert the new event as row 31, pushing the existing rows 31-33 down to 32-34.
$ws1.Rows(31).Insert()

$ws1.Cells.Item(31, 1).Value = 30
$ws1.Cells.Item(31, 1).Style = $ws1.Cells.Item(32, 1).Style

# Column B holds a literal "YYYY-MM-DD" label, not a real date - force text
# storage (Excel would otherwise auto-convert it to a date serial), then
# restore the plain/default cell style copied from a sibling cell.
$ws1.Cells.Item(31, 2).NumberFormat = "@"
$ws1.Cells.Item(31, 2).Value = "2024-06-22"
$ws1.Cells.Item(31, 2).Style = $ws1.Cells.Item(32, 2).Style

$ws1.Cells.Item(31, 3).Value = "广州·第六届淋唔到动漫嘉年华"
$ws1.Cells.Item(31, 4).Value = "沿江东三路15号 广州1978文化创意园"
$ws1.Cells.Item(31, 5).Value = "2024.06.22 10:00-06.23 17:00"
$ws1.Cells.Item(31, 6).Value = 0
$ws1.Cells.Item(31, 7).Value = "不可售"
$ws1.Cells.Item(31, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85554"
$ws1.Cells.Item(31, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/MtLwpx7j1715570717678.jpeg"

# ---------------------------------------------------------------------
# Sheet "演出" (index 2) - single counter bump, no structural change
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(7, 6).Value = 33

# ---------------------------------------------------------------------
# Sheet "全部类型" (index 4) - same counter bumps as sheet 1, plus the
# same newly discovered event inserted as row 35 (pushing 35-41 to 36-42)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$updates4 = @{
    4  = 148
    5  = 1779
    7  = 1061
    8  = 2208
    9  = 2125
    10 = 1111
    11 = 607
    17 = 44
    19 = 97
    22 = 1588
    24 = 727
    26 = 12258
    27 = 12310
    28 = 910
    31 = 38
    33 = 371
    34 = 1923
}

foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}

$ws4.Rows(35).Insert()

$ws4.Cells.Item(35, 1).Value = 34
$ws4.Cells.Item(35, 1).Style = $ws4.Cells.Item(36, 1).Style

$ws4.Cells.Item(35, 2).NumberFormat = "@"
$ws4.Cells.Item(35, 2).Value = "2024-06-22"
$ws4.Cells.Item(35, 2).Style = $ws4.Cells.Item(36, 2).Style

$ws4.Cells.Item(35, 3).Value = "广州·第六届淋唔到动漫嘉年华"
$ws4.Cells.Item(35, 4).Value = "沿江东三路15号 广州1978文化创意园"
$ws4.Cells.Item(35, 5).Value = "2024.06.22 10:00-06.23 17:00"
$ws4.Cells.Item(35, 6).Value = 0
$ws4.Cells.Item(35, 7).Value = "不可售"
$ws4.Cells.Item(35, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85554"
$ws4.Cells.Item(35, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/MtLwpx7j1715570717678.jpeg"
